# Applies the DivisionTemplate.docx edit described by the commit:
#   "Add comments, fix some bugs."
#
# Concretely:
#   1. Move the "_GoBack" bookmark from the very last paragraph of the
#      document to just before the second "Модель " run (the one bookmarked
#      "model2"), which is exactly what Word does automatically when the
#      user's last edit position was there before save - re-adding the
#      "_GoBack" bookmark at a new location removes the old one and Word
#      renumbers every w:id sequentially.
#   2. Bump the font size (both w:sz and w:szCs) of the last table row
#      (9 cells) in the first table from 8pt (16 half-points) to 10pt
#      (20 half-points).
#   3. Reduce the top page margin from 56.7pt (1134 twips) to 42.55pt
#      (851 twips).

$d = $word.ActiveDocument

# --- 1. Re-anchor the "_GoBack" bookmark -----------------------------------
# Find the 2nd occurrence of "Модель " (the first belongs to the "model"
# bookmark near the top of the doc; the second is immediately followed by
# the "model2" bookmark).
$r = $d.Content
[void]$r.Find.Execute("Модель ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(0)
[void]$r.Find.Execute("Модель ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.Collapse(1)
[void]$d.Bookmarks.Add("_GoBack", $r)

# --- 2. Grow the font size used by the last row of the first table --------
$tbl = $d.Tables.Item(1)
$lastRow = $tbl.Rows.Item($tbl.Rows.Count)
$lastRow.Range.Font.Size = 10
$lastRow.Range.Font.SizeBi = 10

# --- 3. Shrink the top page margin -----------------------------------------
$d.PageSetup.TopMargin = 42.55
